$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value2 = "nourhan.mahmoud@med.asu.edu.eg , shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G9").Value2 = "wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G14").Value2 = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G15").Value2 = "Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G27").Value2 = "System, randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G29").Value2 = "ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G34").Value2 = "nourhan.mahmoud@med.asu.edu.eg , shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G40").Value2 = "wessam.atef@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G45").Value2 = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G46").Value2 = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G58").Value2 = "System, randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G60").Value2 = "ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G64").Value2 = "asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G65").Value2 = "Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G71").Value2 = "Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
$ws.Range("G76").Value2 = "Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G77").Value2 = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G91").Value2 = "ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G95").Value2 = "mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G96").Value2 = "alshimaa.atef@med.asu.edu.egm, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg "
$ws.Range("G102").Value2 = "aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G107").Value2 = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G108").Value2 = "Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G120").Value2 = "nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G122").Value2 = "eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G126").Value2 = "mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G127").Value2 = "alshimaa.atef@med.asu.edu.egm, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg "
$ws.Range("G133").Value2 = "aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G138").Value2 = "Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G139").Value2 = "Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G140").Value2 = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G151").Value2 = "user@user.com, Amr-Saeed@med.asu.edu.eg`n, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg"
$ws.Range("G153").Value2 = "eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G157").Value2 = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G158").Value2 = "Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G164").Value2 = "Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
$ws.Range("G169").Value2 = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G170").Value2 = "Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G182").Value2 = "nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G184").Value2 = "ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
